$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing headers
# (copy format from H1, the neighboring header cell, then overwrite the text)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data values for I2:J70 (columns I0 / IF), one @(I, J) pair per row starting at row 2
$data = @(
    @(7,8),   @(8,8),   @(7,7),   @(8,8),   @(8,8),   @(8,8),   @(8,8),   @(1,1),
    @(8,8),   @(7,8),   @(7,7),   @(7,7),   @(7,7),   @(8,8),   @(7,7),   @(8,8),
    @(8,8),   @(8,8),   @(7,7),   @(7,7),   @(10,10), @(9,9),   @(7,7),   @(10,10),
    @(6,7),   @(9,9),   @(8,8),   @(8,8),   @(1,1),   @(8,8),   @(6,6),   @(1,1),
    @(8,8),   @(7,7),   @(7,7),   @(7,7),   @(7,7),   @(8,8),   @(7,7),   @(7,7),
    @(1,1),   @(8,8),   @(9,9),   @(8,8),   @(9,9),   @(8,8),   @(7,7),   @(8,8),
    @(6,6),   @(9,9),   @(8,8),   @(8,8),   @(1,1),   @(7,7),   @(9,9),   @(8,8),
    @(9,9),   @(8,8),   @(8,8),   @(9,9),   @(7,8),   @(9,9),   @(1,1),   @(7,7),
    @(8,8),   @(1,1),   @(8,8),   @(9,9),   @(5,5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
